$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight two "best" values in the original table in bold ---
$ws.Range("I5").Font.Bold = $true
$ws.Range("G17").Font.Bold = $true

# --- New section header (row 28), merged + centered ---
$ws.Range("F28:H28").Merge()
$ws.Range("F28").Value = "Spearman Rank Correlation for 2 Datasets"
$ws.Range("F28:H28").HorizontalAlignment = -4108

# --- New column headers (row 29) ---
$ws.Range("F29").Value = "Score Metric"
$ws.Range("G29").Value = "Summaries"
$ws.Range("H29").Value = "Essays"

# --- New data table (rows 30-44): Score Metric / Summaries / Essays ---
$ws.Range("F30").Value = "BLEU"
$ws.Range("G30").Value = 0.0367668094941715
$ws.Range("G30").NumberFormat = "0.000"
$ws.Range("H30").Value = 0.61538424199999997
$ws.Range("H30").NumberFormat = "0.000"
$ws.Range("H30").Font.Bold = $true
$ws.Range("J30").NumberFormat = "0.00E+00"

$ws.Range("F31").Value = "ROUGE-L"
$ws.Range("G31").Value = 0.102188827507858
$ws.Range("G31").NumberFormat = "0.000"
$ws.Range("H31").Value = 0.33663348700000001
$ws.Range("H31").NumberFormat = "0.000"
$ws.Range("J31").NumberFormat = "0.00E+00"

$ws.Range("F32").Value = "glove_wms"
$ws.Range("G32").Value = 0.180172319775732
$ws.Range("G32").NumberFormat = "0.000"
$ws.Range("H32").Value = 0.42856114499999998
$ws.Range("H32").NumberFormat = "0.000"
$ws.Range("J32").NumberFormat = "0.00E+00"

$ws.Range("F33").Value = "glove_sms"
$ws.Range("G33").Value = 0.25732861378390298
$ws.Range("G33").NumberFormat = "0.000"
$ws.Range("H33").Value = 0.44947218799999999
$ws.Range("H33").NumberFormat = "0.000"
$ws.Range("J33").NumberFormat = "0.00E+00"

$ws.Range("F34").Value = "glove_s+wms"
$ws.Range("G34").Value = 0.21389133917186001
$ws.Range("G34").NumberFormat = "0.000"
$ws.Range("H34").Value = 0.48756354800000001
$ws.Range("H34").NumberFormat = "0.000"
$ws.Range("J34").NumberFormat = "0.00E+00"

$ws.Range("F35").Value = "elmo_wms"
$ws.Range("G35").Value = 0.16030232756682
$ws.Range("G35").NumberFormat = "0.000"
$ws.Range("H35").Value = 0.44036108400000001
$ws.Range("H35").NumberFormat = "0.000"
$ws.Range("J35").NumberFormat = "0.00E+00"

$ws.Range("F36").Value = "elmo_sms"
$ws.Range("G36").Value = 0.25307399956494098
$ws.Range("G36").NumberFormat = "0.000"
$ws.Range("H36").Value = 0.43804670200000001
$ws.Range("H36").NumberFormat = "0.000"
$ws.Range("J36").NumberFormat = "0.00E+00"

$ws.Range("F37").Value = "elmo_s+wms"
$ws.Range("G37").Value = 0.20347261877414299
$ws.Range("G37").NumberFormat = "0.000"
$ws.Range("H37").Value = 0.48320539099999998
$ws.Range("H37").NumberFormat = "0.000"
$ws.Range("J37").NumberFormat = "0.00E+00"

$ws.Range("F38").Value = "bert_wms"
$ws.Range("G38").Value = 0.16889596200715501
$ws.Range("G38").NumberFormat = "0.000"
$ws.Range("H38").Value = 0.36768276
$ws.Range("H38").NumberFormat = "0.000"
$ws.Range("J38").NumberFormat = "0.00E+00"

$ws.Range("F39").Value = "bert_sms"
$ws.Range("G39").Value = 0.22895971006059601
$ws.Range("G39").NumberFormat = "0.000"
$ws.Range("H39").Value = 0.41175308500000002
$ws.Range("H39").NumberFormat = "0.000"
$ws.Range("J39").NumberFormat = "0.00E+00"

$ws.Range("F40").Value = "bert_s+wms"
$ws.Range("G40").Value = 0.19876171718671101
$ws.Range("G40").NumberFormat = "0.000"
$ws.Range("H40").Value = 0.413584646
$ws.Range("H40").NumberFormat = "0.000"
$ws.Range("J40").NumberFormat = "0.00E+00"

$ws.Range("F41").Value = "BERTScore_P"
$ws.Range("G41").Value = 0.178426557580073
$ws.Range("G41").NumberFormat = "0.000"
$ws.Range("H41").Value = -0.13474971799999999
$ws.Range("H41").NumberFormat = "0.000"
$ws.Range("J41").NumberFormat = "0.00E+00"

$ws.Range("F42").Value = "BERTScore_R"
$ws.Range("G42").Value = 0.26315372122443498
$ws.Range("G42").NumberFormat = "0.000"
$ws.Range("G42").Font.Bold = $true
$ws.Range("H42").Value = 0.53577657199999995
$ws.Range("H42").NumberFormat = "0.000"
$ws.Range("J42").NumberFormat = "0.00E+00"

$ws.Range("F43").Value = "BERTScore_F1"
$ws.Range("G43").Value = 0.25525405987705801
$ws.Range("G43").NumberFormat = "0.000"
$ws.Range("H43").Value = 0.147088571
$ws.Range("H43").NumberFormat = "0.000"
$ws.Range("J43").NumberFormat = "0.00E+00"

$ws.Range("F44").Value = "SentBERT"
$ws.Range("G44").Value = 0.11540827300000001
$ws.Range("G44").NumberFormat = "0.000"
$ws.Range("H44").Value = 0.29691668700000001
$ws.Range("H44").NumberFormat = "0.000"
$ws.Range("J44").NumberFormat = "0.00E+00"

# --- Column widths: narrower G, new custom-width H, I keeps its autofit width ---
$ws.Columns("G").ColumnWidth = 12.67
$ws.Columns("H").ColumnWidth = 8

# --- View: zoom in on the sheet and move the active selection ---
$win = $excel.ActiveWindow
$win.Zoom = 130
$ws.Range("I16").Select()
